$d = $word.ActiveDocument
$d.TrackRevisions = $false

# Change 1: "MA_G08_01_CO" -> "MA_08_01_CO"
$d.Content.Find.Execute("MA_G08_01_CO", $true, $false, $false, $false, $false,
                         $true, 1, $false, "MA_08_01_CO", 1)

# Change 2: merge the two runs "ejercicios:" and " " into a single run
# containing "ejercicios: ", without touching the preceding
# "Resuelve los siguientes " run.
$findRng = $d.Content
$findRng.Find.Execute("ejercicios:")
$spaceStart = $findRng.End
$spaceRng = $d.Range($spaceStart, $spaceStart + 1)
$spaceRng.Find.Execute(" ", $true, $false, $false, $false, $false, $true, 1,
                        $false, " ", 1)
